# Apply "corrections in the vocabulary" edit to the workbook.
#
# Summary of changes (per the target diff):
#  - On sheet "factor_binary": the contents (value + formatting) of
#    column B and column C were swapped on rows 5, 6 and 20 (the
#    Yes/No, Correct/Incorrect and well-informed/not-well-informed
#    pairs had been entered in the wrong columns).
#  - The "factor_binary" sheet becomes the active/selected sheet, with
#    cell B6 selected (previously D3 was selected there).
#  - The "keep_numeric_names" sheet (previously the active sheet) is
#    no longer the selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("factor_binary")

# Swaps the B/C cell pair (value + number format/font/style) on a given row.
function Swap-BC-Row($sheet, $row) {
    $staging = $sheet.Range("Z100")
    $sheet.Cells.Item($row, 2).Copy($staging)
    $sheet.Cells.Item($row, 3).Copy($sheet.Cells.Item($row, 2))
    $staging.Copy($sheet.Cells.Item($row, 3))
    $staging.Clear()
}

Swap-BC-Row $ws 5
Swap-BC-Row $ws 6
Swap-BC-Row $ws 20

# Make factor_binary the active sheet (this also clears tabSelected on
# whichever sheet was previously active, i.e. keep_numeric_names) and
# select cell B6.
$ws.Activate()
$ws.Range("B6").Select()
